$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# 1. Schweilarbeiten -> Schweiarbeiten
Replace-Text "Bei Schneid- und Schweilarbeiten in feuer- und explosionsgefährdeten Räumen müssen vor Beginn der Arbeiten alle Gefahrengüter entfernt werden." `
             "Bei Schneid- und Schweiarbeiten in feuer- und explosionsgefährdeten Räumen müssen vor Beginn der Arbeiten alle Gefahrengüter entfernt werden."

# 2. Schweijßgerät -> Schweißgerät
Replace-Text "Was muss bei Wartungsarbeiten an einmem Schweijßgerät immer beachtet werden?" `
             "Was muss bei Wartungsarbeiten an einmem Schweißgerät immer beachtet werden?"

# 3. Latenhebemagnet -> Lastenhebemagnet
Replace-Text "Dürfen die Gasflaschen mit Hilfe von Latenhebemagnet transportiert werden?" `
             "Dürfen die Gasflaschen mit Hilfe von Lastenhebemagnet transportiert werden?"

# 4. "wird odie Schweißnaht ... gegn Luftzutritt" -> "wird die Schweißnaht ... gegen Luftzutritt"
Replace-Text "Beim Schutzgasschweißen wird odie Schweißnaht durch einen Schutzgasmantel gegn Luftzutritt geschützt." `
             "Beim Schutzgasschweißen wird die Schweißnaht durch einen Schutzgasmantel gegen Luftzutritt geschützt."

# 5. geigenet -> geeignet
Replace-Text "Was ist für die Gütesicherung der Schweißnaht nicht geigenet?" `
             "Was ist für die Gütesicherung der Schweißnaht nicht geeignet?"

# 6. Schweilßfolgeplan -> Schweißfolgeplan (3 occurrences)
Replace-Text "Schweilßfolgeplan" "Schweißfolgeplan"

# 7. MAG- Schweilßen -> MAG- Schweißen
Replace-Text "Durch richtige Parametereinstellung kann man die Spritzerbildung beim MAG- Schweilßen verringern." `
             "Durch richtige Parametereinstellung kann man die Spritzerbildung beim MAG- Schweißen verringern."

# 8. Härteneignung -> Härteneigung (2 occurrences)
Replace-Text "Härteneignung" "Härteneigung"

# 9. Welche auswirkung -> Welche Auswirkung
Replace-Text "Welche auswirkung hat Kohlenstoff beim unsachgemäßen Schweißen?" `
             "Welche Auswirkung hat Kohlenstoff beim unsachgemäßen Schweißen?"

# 10. MAG-Schweilßen verwendet -> MAG-Schweißen verwendet
Replace-Text "CO₂ kann zum MAG-Schweilßen verwendet werden." `
             "CO₂ kann zum MAG-Schweißen verwendet werden."

# 11. MIG-Schweilßen -> MIG-Schweißen (within the MAG/MIG explanation sentence)
Replace-Text "Zwischen MAG- und MIG-Verfahren gibt es keine generelle Austauschbarkeit der Gase. MAG-Schweißen --> Aktives Gas MIG-Schweilßen --> Inertes Gas" `
             "Zwischen MAG- und MIG-Verfahren gibt es keine generelle Austauschbarkeit der Gase. MAG-Schweißen --> Aktives Gas MIG-Schweißen --> Inertes Gas"
